$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Groom PBI014 (row 9): refine the "labels" story with the new search-results idea
$ws.Range("B9").Value = "As an administrator`nI want every saved exercise and its questions to have labels of question-type, topic and anything else that can enhance search results`nBecause the user wants auto-generated tests"
$ws.Range("E9").Value = "Groomed"
$ws.Rows.Item(9).RowHeight = 60

# Groom PBI017 (row 11): elaborate on labels/folders syncing to the database
$ws.Range("B11").Value = "As a user`nI don't want labels to be fixed`nBecause I might think of topics of my own`n(distinguish between category and (optional) label? Users use labels/folders to order their tests; the names of which are copied to the database for better search results)"
$ws.Range("E11").Value = "Groomed"
$ws.Rows.Item(11).RowHeight = 75

# Groom PBI015 (row 10): mention copying folder names to database labels
$ws.Range("B10").Value = "As a user`nI don't want to add a label to every question (require topic? Copy folder names etc to database labels)`nBecause I want fast"
$ws.Range("E10").Value = "Groomed"
$ws.Rows.Item(10).RowHeight = 60

# Leave the cursor/viewport near the freshly groomed rows
$ws.Range("B9").Select()
